$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "floatmod"
$ws.Range("B15").Value = 500
$ws.Range("C15").Value = "浮力模块"

$ws.Range("C16").Select()
